$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 98.912777
$ws.Range("H2").Value = 296.738331
$ws.Range("I2").Value = 0.8120825131376513
$ws.Range("J2").Value = 0.8120825131376513
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.7065936666666666
$ws.Range("N2").Value = 2.119781
$ws.Range("O2").Value = 0.005187843618793344
$ws.Range("P2").Value = 0.005187843618793344
$ws.Range("Q2").Value = 69.89114178061233
$ws.Range("R2").Value = 629.020276025511
$ws.Range("S2").Value = 0.004212957083714827
$ws.Range("T2").Value = 0.004212957083714827

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 98.912777
$ws.Range("H3").Value = 296.738331
$ws.Range("I3").Value = 0.8120825131376513
$ws.Range("J3").Value = 0.8120825131376513
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 111.9320066666667
$ws.Range("N3").Value = 335.79602
$ws.Range("O3").Value = 0.8218100075305903
$ws.Range("P3").Value = 0.8218100075305903
$ws.Range("Q3").Value = 11071.50561458251
$ws.Range("R3").Value = 99643.55053124262
$ws.Range("S3").Value = 0.6673775362371139
$ws.Range("T3").Value = 0.6673775362371139

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 98.912777
$ws.Range("H4").Value = 296.738331
$ws.Range("I4").Value = 0.8120825131376513
$ws.Range("J4").Value = 0.8120825131376513
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 23.563205
$ws.Range("N4").Value = 70.689615
$ws.Range("O4").Value = 0.1730021488506163
$ws.Range("P4").Value = 0.1730021488506163
$ws.Range("Q4").Value = 2330.702041570285
$ws.Range("R4").Value = 20976.31837413257
$ws.Range("S4").Value = 0.1404920198168225
$ws.Range("T4").Value = 0.1404920198168225

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 17.04862266666667
$ws.Range("H5").Value = 51.14586800000001
$ws.Range("I5").Value = 0.1399706767982279
$ws.Range("J5").Value = 0.1399706767982279
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7065936666666666
$ws.Range("N5").Value = 2.119781
$ws.Range("O5").Value = 0.005187843618793344
$ws.Range("P5").Value = 0.005187843618793344
$ws.Range("Q5").Value = 12.04644880165645
$ws.Range("R5").Value = 108.418039214908
$ws.Range("S5").Value = 0.0007261459824458724
$ws.Range("T5").Value = 0.0007261459824458724

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 17.04862266666667
$ws.Range("H6").Value = 51.14586800000001
$ws.Range("I6").Value = 0.1399706767982279
$ws.Range("J6").Value = 0.1399706767982279
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 111.9320066666667
$ws.Range("N6").Value = 335.79602
$ws.Range("O6").Value = 0.8218100075305903
$ws.Range("P6").Value = 0.8218100075305903
$ws.Range("Q6").Value = 1908.286545982818
$ws.Range("R6").Value = 17174.57891384536
$ws.Range("S6").Value = 0.1150293029536135
$ws.Range("T6").Value = 0.1150293029536135

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 17.04862266666667
$ws.Range("H7").Value = 51.14586800000001
$ws.Range("I7").Value = 0.1399706767982279
$ws.Range("J7").Value = 0.1399706767982279
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 23.563205
$ws.Range("N7").Value = 70.689615
$ws.Range("O7").Value = 0.1730021488506163
$ws.Range("P7").Value = 0.1730021488506163
$ws.Range("Q7").Value = 401.7201908623134
$ws.Range("R7").Value = 3615.481717760821
$ws.Range("S7").Value = 0.02421522786216854
$ws.Range("T7").Value = 0.02421522786216854

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.839988000000001
$ws.Range("H8").Value = 17.519964
$ws.Range("I8").Value = 0.0479468100641207
$ws.Range("J8").Value = 0.04794681006412069
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.7065936666666666
$ws.Range("N8").Value = 2.119781
$ws.Range("O8").Value = 0.005187843618793344
$ws.Range("P8").Value = 0.005187843618793344
$ws.Range("Q8").Value = 4.126498534209333
$ws.Range("R8").Value = 37.138486807884
$ws.Range("S8").Value = 0.0002487405526326451
$ws.Range("T8").Value = 0.000248740552632645

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.839988000000001
$ws.Range("H9").Value = 17.519964
$ws.Range("I9").Value = 0.0479468100641207
$ws.Range("J9").Value = 0.04794681006412069
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 111.9320066666667
$ws.Range("N9").Value = 335.79602
$ws.Range("O9").Value = 0.8218100075305903
$ws.Range("P9").Value = 0.8218100075305903
$ws.Range("Q9").Value = 653.6815757492534
$ws.Range("R9").Value = 5883.13418174328
$ws.Range("S9").Value = 0.03940316833986281
$ws.Range("T9").Value = 0.03940316833986281

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 5.839988000000001
$ws.Range("H10").Value = 17.519964
$ws.Range("I10").Value = 0.0479468100641207
$ws.Range("J10").Value = 0.04794681006412069
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 23.563205
$ws.Range("N10").Value = 70.689615
$ws.Range("O10").Value = 0.1730021488506163
$ws.Range("P10").Value = 0.1730021488506163
$ws.Range("Q10").Value = 137.60883444154
$ws.Range("R10").Value = 1238.47950997386
$ws.Range("S10").Value = 0.008294901171625237
$ws.Range("T10").Value = 0.008294901171625237

